# repeat_loans.xlsx - "update of results using corrected definition of apr"
#
# The sheet's data cells are formula cells pulling cached values from a
# closed external workbook (=[1]repeat_loans!B5, etc. - a link to
# reg_results/repeat_loans.csv). The commit re-ran the regressions with a
# corrected APR definition, which changed the *cached* numbers that were
# displayed for a number of coefficients / SEs / Ns / R-squared / control
# means in the table, while leaving every other cell (headers, labels,
# unaffected columns) untouched.
#
# Each assignment below uses a leading apostrophe so Excel stores the text
# exactly as printed in the source table (e.g. "0.30", "(0.035)", "0.067*")
# instead of re-interpreting it as a number and dropping trailing zeros /
# introducing floating point noise.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - "Forced commitment" coefficients
$ws.Range("B5").Value = "'0.067*"
$ws.Range("D5").Value = "'-0.0046"
$ws.Range("E5").Value = "'0.036"
$ws.Range("F5").Value = "'0.037***"
$ws.Range("I5").Value = "'7.72**"

# Row 6 - Row 5's standard errors
$ws.Range("B6").Value = "'(0.035)"
$ws.Range("C6").Value = "'(0.047)"
$ws.Range("D6").Value = "'(0.035)"
$ws.Range("E6").Value = "'(0.030)"
$ws.Range("F6").Value = "'(0.013)"
$ws.Range("G6").Value = "'(0.027)"
$ws.Range("I6").Value = "'(3.07)"

# Row 7 - "Choice commitment" coefficients
$ws.Range("B7").Value = "'0.040"
$ws.Range("C7").Value = "'0.051"
$ws.Range("D7").Value = "'0.026"
$ws.Range("E7").Value = "'0.027"
$ws.Range("F7").Value = "'0.0098"
$ws.Range("I7").Value = "'1.72"

# Row 8 - Row 7's standard errors
$ws.Range("B8").Value = "'(0.031)"
$ws.Range("C8").Value = "'(0.042)"
$ws.Range("D8").Value = "'(0.034)"
$ws.Range("E8").Value = "'(0.027)"
$ws.Range("F8").Value = "'(0.0087)"
$ws.Range("G8").Value = "'(0.026)"
$ws.Range("I8").Value = "'(2.59)"

# Row 10 - Observations
$ws.Range("C10").Value = "'2168"
$ws.Range("D10").Value = "'2273"
$ws.Range("I10").Value = "'1577"

# Row 11 - R-squared
$ws.Range("B11").Value = "'0.003"
$ws.Range("F11").Value = "'0.006"
$ws.Range("I11").Value = "'0.011"

# Row 12 - Control Mean
$ws.Range("B12").Value = "'0.32"
$ws.Range("C12").Value = "'0.36"
$ws.Range("D12").Value = "'0.29"
$ws.Range("E12").Value = "'0.28"
$ws.Range("F12").Value = "'0.020"
$ws.Range("I12").Value = "'32.9"
